$d = $word.ActiveDocument

# --- Locate the insertion point -------------------------------------------------
# The 5 new bullet paragraphs are inserted right before the paragraph that
# starts with "Nhn 입사 결정..." (the "N" / "hn" run pair), which currently
# sits right after "플레이 스타일도 ...".
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Nhn")) {
        $target = $p
        break
    }
}

$insPoint = $target.Range.Start
$ins = $d.Range($insPoint, $insPoint)

# Build the five new list paragraphs in one shot so they inherit the
# "List Paragraph" (a5) style / numId=1 / ilvl=0 formatting of the
# paragraph that follows them (the "Nhn..." paragraph already uses that
# same style+level).
$cr = [char]13
$newText = "해상도는 1280 * 720으로 결정" + $cr + `
           "UI를 유니티 4.6에서 지원하는 기본UI로 갈 것인지 NGUI로 갈 것인지는 추후 판단하기로 보류" + $cr + `
           "매수 수요일 정기 모임을 가지기로 결정." + $cr + `
           "장소 : " + $cr + `
           "시간 : " + $cr
$ins.InsertBefore($newText)

# --- Re-find the two new paragraphs that must sit at list level 2 (ilvl=1) ------
$p1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("장소")) {
        $p1 = $p
        break
    }
}
$p1.Range.ListFormat.ListLevelNumber = 2

$p2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("시간")) {
        $p2 = $p
        break
    }
}
$p2.Range.ListFormat.ListLevelNumber = 2

# --- Move the hidden "_GoBack" bookmark to the end of the last paragraph typed --
# ("시간 : ") to mirror where Word would leave it after the edit.
$goBackRange = $d.Range($p2.Range.End - 1, $p2.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)
